$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Z2").Value = 12
$ws.Range("AE2").Value = 21
$ws.Range("AH2").Value = 11
$ws.Range("AJ2").Value = 19
$ws.Range("AM2").Value = 51
$ws.Range("AO2").Value = 9
$ws.Range("AP2").Value = 23
$ws.Range("AR2").Value = 51
$ws.Range("AW2").Value = 7
$ws.Range("BB2").Value = 451

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("I4").Value = 2.15
$ws.Range("J4").Value = 4.33
$ws.Range("L4").Value = 3
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Y4").Value = 13
$ws.Range("AI4").Value = 9
$ws.Range("AJ4").Value = 10
$ws.Range("AO4").Value = 21
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 29

# Row 5
$ws.Range("G5").Value = 1.53
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 8
$ws.Range("Y5").Value = 9.5
$ws.Range("AH5").Value = 12
$ws.Range("AO5").Value = 8
$ws.Range("AP5").Value = 26

# Row 6
$ws.Range("Q6").Value = 2.35
$ws.Range("R6").Value = 1.57
